$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.311.85"

$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "1.650.54"

$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.75"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0631"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.06"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  -0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").Value = "1.881.47"

$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("D13").Value = "1.655.91"

$ws.Range("E13").Value = "  -0.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  -1.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.542"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +1.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.91"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +1.25%  "

$ws.Range("D17").Value = "27.298.76"

$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "0.0₃0741"

$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.95"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  -1.15%  "

$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +1.48%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("E24").Value = "  -0.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.84"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.56"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.86"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.03"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +0.57%  "

$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").Value = "1.258.58"

$ws.Range("E35").Value = "  -0.54%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.544"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  +1.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.843"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("E41").Value = "  -0.99%  "

$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("E43").Value = "  +4.63%  "

$ws.Range("D44").Value = "1.791.60"

$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.34"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.93"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("B48").Value = "BabyDogeCoin"

$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

$ws.Range("D48").Value = "0.0₆0100"

$ws.Range("E48").Value = "  +15.69%  "

$ws.Range("B49").Value = "Cronos"

$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -0.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("E51").Value = "  -0.94%  "
